$d = $word.ActiveDocument

# The 3rd table on the page is the "manager / job" signer table
# (columns: numbering | name | job). It currently auto-sizes; the
# update pins it to a fixed width and widens the 3rd ("job") column.
$t = $d.Tables.Item(3)

# Switch the table from "auto" width to a fixed width of 8651 twips
# (432.55 pt) -- wdPreferredWidthPoints-style fixed layout.
$t.PreferredWidthType = 3
$t.PreferredWidth = 432.55

# Widen the 3rd column from 3645 twips (182.25 pt) to 3933 twips
# (196.65 pt); the other two columns (425 / 4293 twips) stay as-is.
$col = $t.Columns.Item(3)
$col.PreferredWidth = 196.65
